$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (events that were removed in this update), shifting rows up
$ws.Rows("4:5").Delete()

# Force the numeric-looking ID columns (D,E,F) in rows 2-3 to be stored as text,
# matching the original inlineStr convention used throughout this sheet.
$idCols = $ws.Range("D2:F3")
$idCols.NumberFormat = "@"

# Row 2: new event - No Seat Belt / LUIS IBARRA
$ws.Range("A2").Value = "281474991205821-1749752036273"
$ws.Range("B2").Value = "No Seat Belt"
$ws.Range("C2").Value = "2025-06-12T12:13:56.273"
$ws.Range("D2").Value = "281474991205821"
$ws.Range("E2").Value = "148"
$ws.Range("F2").Value = "51834015"
$ws.Range("G2").Value = "LUIS IBARRA"
$ws.Range("H2").Value = 20.56618713
$ws.Range("I2").Value = -103.45906243
$ws.Range("J2").Value = 0

$k2bytes = [System.Convert]::FromBase64String("aHR0cHM6Ly9zMy5zYW1zYXJhLmNvbS9zYW1zYXJhLWRhc2hjYW0tdmlkZW9zLzQwMDYxMjQvMjgxNDc0OTkxMjA1ODIxLzE3NDk3NTIwMzM3NzMvRk5DZGxsSzBrUS1jYW1lcmEtdmlkZW8tc2VnbWVudC1kcml2ZXItMTc0OTc1MjAzNjI3My5tcDQ/WC1BbXotQWxnb3JpdGhtPUFXUzQtSE1BQy1TSEEyNTYmWC1BbXotQ3JlZGVudGlhbD1BU0lBM0xZM1JOV1NORUNORzVNTCUyRjIwMjUwNjEzJTJGdXMtd2VzdC0yJTJGczMlMkZhd3M0X3JlcXVlc3QmWC1BbXotRGF0ZT0yMDI1MDYxM1QxNjAwMzdaJlgtQW16LUV4cGlyZXM9Mjg4MDAmWC1BbXotU2VjdXJpdHktVG9rZW49SVFvSmIzSnBaMmx1WDJWakVDNGFDWFZ6TFhkbGMzUXRNaUpJTUVZQ0lRRDdsYjh2RmJVUmJSViUyQmtUN1ZMYVE5QlBYZDZUbGp6SElIcUd3dG9lUkE4d0loQUlvMFJWVXptWmFRdGNxOXNTdlBTTGNVaWp4JTJCWmhoZkJGN2I1c1RUT2JEVkt0MERDQmNRQkJvTU56Z3hNakEwT1RReU1qUTBJZ3c1eHZOeUZDRjU4bVMlMkJTQm9xdWdNVU1MZEZkVW15cHJTQSUyRkN4YWJudTlmV2xJY3doNTl2SUhxZzA2ckZsN050Y1pRYlhya0dQQVp2ZksyQlZHWHdEdUxVd1lxVWtzYnZSek5iODVVdVliSSUyQkJhNkw0dmpwbUkwSUx4emtOMU5OTGJ2VTNWQkdTeFFzJTJCaGF2aVJ0TDZ2SlhaaEs2Tnk2TXU1RWlMWjhXWmNhR0Rvdk8xb2lJYWFWeFhvM1JBU1RRckd3NWR4S0V1QkFoaHUlMkI4ckxKSSUyRmN0WEwlMkZQdGM1bUtWJTJCNllGY2J1QWZoSTUyemgzSVdDUEx2dVEzQk4lMkY5eXJoVWxBMDgyQTZHcGxyRGYyVEJBTCUyRmxib3NjeXYyT3hpUkdUeUdtNk5odFBhTmFRSUl2UktjeWZzd2Z4c2tOcHhVcWdRVXEyWklkaTFuMVlNZkM5aEZJajVHNUJCenZHVTNqd0VmTVlReTN1anNIeklVVFlxc2lzMjhZSiUyRmolMkJYTWpIbW54SW52eVZ1amJteW5DdEhNemZhTXJ4R240MFNOQVdPUE16TGRpVGN3cEtwdFolMkJLTUNjUzklMkZpQ3clMkYxWHRRT2E5V0I2TGpHMW5SSDBiNTFrVzNVb0w5Skl1TEVUSWtxQ3BhbiUyQkdldjY5Wkc0OCUyRmo1WXZ5QTNVJTJGVE55czBvTFQ4eVN0Nk9WZU5XbXA0ZDJkdGhGZHUwU1V3VXdDUERUVkRXbnY0dk9jVWJHemFuSGM2QjFBc29xMW1iNHB6NTI1R09tY1YwSHF1bVl3UEtyaDgzQWZvMzhBTXQlMkJEN1ZodE1QUGdzTUlHT3FRQk5qT05HN3I0bndJZ0olMkZBZEI5d1pPOVoyMnRMenE0QU1OJTJCJTJCSWVQSDIxNlVreFNZOENicXVrWjl0R2FaeGlURXJXNEZaOGNRajJ6UXNUUHZ2NVNoTWNsYmVnbFRzRnNjenpGeTB4T1QlMkJsVkxzZjhaMW05WURWWTRhSDFWOWlsMlpPQkoyZFdMY0xCeXRVRmJvVXVpY2pYZ3puSGJaYkJrZWFqZ0lvYkQlMkJYQ1Z1OGF6em44b29pJTJGdXE3RmdKbDdMJTJCUXpkczl0bGx0NnFpdmVYdGNIdWFQbktWSjFBJTNEJlgtQW16LVNpZ25lZEhlYWRlcnM9aG9zdCZyZXNwb25zZS1leHBpcmVzPVNhdCUyQyUyMDE0JTIwSnVuJTIwMjAyNSUyMDAwJTNBMDAlM0EzNyUyMEdNVCZYLUFtei1TaWduYXR1cmU9NTY5MTg5YzM0MWJiOTM5YjJkMGI3ZDEyNGRlMzJiZjc4YTM3ZTE2YTNhYWIxZTgzYzY3ZTdkNGQ0Njg5N2E4NQ==")
$ws.Range("K2").Value = [System.Text.Encoding]::UTF8.GetString($k2bytes)

$l2bytes = [System.Convert]::FromBase64String("aHR0cHM6Ly9zMy5zYW1zYXJhLmNvbS9zYW1zYXJhLWN2ZGF0YS80MDA2MTI0LzI4MTQ3NDk5MTIwNTgyMS8xNzQ5NzUyMDMzNzczL2YxTmZlT2NqY20tY2FtZXJhLXZpZGVvLXNlZ21lbnQtMTc0OTc1MjAzNjI3My5hdWRpby5tcDQ/WC1BbXotQWxnb3JpdGhtPUFXUzQtSE1BQy1TSEEyNTYmWC1BbXotQ3JlZGVudGlhbD1BU0lBM0xZM1JOV1NORUNORzVNTCUyRjIwMjUwNjEzJTJGdXMtd2VzdC0yJTJGczMlMkZhd3M0X3JlcXVlc3QmWC1BbXotRGF0ZT0yMDI1MDYxM1QxNjAwMzdaJlgtQW16LUV4cGlyZXM9Mjg4MDAmWC1BbXotU2VjdXJpdHktVG9rZW49SVFvSmIzSnBaMmx1WDJWakVDNGFDWFZ6TFhkbGMzUXRNaUpJTUVZQ0lRRDdsYjh2RmJVUmJSViUyQmtUN1ZMYVE5QlBYZDZUbGp6SElIcUd3dG9lUkE4d0loQUlvMFJWVXptWmFRdGNxOXNTdlBTTGNVaWp4JTJCWmhoZkJGN2I1c1RUT2JEVkt0MERDQmNRQkJvTU56Z3hNakEwT1RReU1qUTBJZ3c1eHZOeUZDRjU4bVMlMkJTQm9xdWdNVU1MZEZkVW15cHJTQSUyRkN4YWJudTlmV2xJY3doNTl2SUhxZzA2ckZsN050Y1pRYlhya0dQQVp2ZksyQlZHWHdEdUxVd1lxVWtzYnZSek5iODVVdVliSSUyQkJhNkw0dmpwbUkwSUx4emtOMU5OTGJ2VTNWQkdTeFFzJTJCaGF2aVJ0TDZ2SlhaaEs2Tnk2TXU1RWlMWjhXWmNhR0Rvdk8xb2lJYWFWeFhvM1JBU1RRckd3NWR4S0V1QkFoaHUlMkI4ckxKSSUyRmN0WEwlMkZQdGM1bUtWJTJCNllGY2J1QWZoSTUyemgzSVdDUEx2dVEzQk4lMkY5eXJoVWxBMDgyQTZHcGxyRGYyVEJBTCUyRmxib3NjeXYyT3hpUkdUeUdtNk5odFBhTmFRSUl2UktjeWZzd2Z4c2tOcHhVcWdRVXEyWklkaTFuMVlNZkM5aEZJajVHNUJCenZHVTNqd0VmTVlReTN1anNIeklVVFlxc2lzMjhZSiUyRmolMkJYTWpIbW54SW52eVZ1amJteW5DdEhNemZhTXJ4R240MFNOQVdPUE16TGRpVGN3cEtwdFolMkJLTUNjUzklMkZpQ3clMkYxWHRRT2E5V0I2TGpHMW5SSDBiNTFrVzNVb0w5Skl1TEVUSWtxQ3BhbiUyQkdldjY5Wkc0OCUyRmo1WXZ5QTNVJTJGVE55czBvTFQ4eVN0Nk9WZU5XbXA0ZDJkdGhGZHUwU1V3VXdDUERUVkRXbnY0dk9jVWJHemFuSGM2QjFBc29xMW1iNHB6NTI1R09tY1YwSHF1bVl3UEtyaDgzQWZvMzhBTXQlMkJEN1ZodE1QUGdzTUlHT3FRQk5qT05HN3I0bndJZ0olMkZBZEI5d1pPOVoyMnRMenE0QU1OJTJCJTJCSWVQSDIxNlVreFNZOENicXVrWjl0R2FaeGlURXJXNEZaOGNRajJ6UXNUUHZ2NVNoTWNsYmVnbFRzRnNjenpGeTB4T1QlMkJsVkxzZjhaMW05WURWWTRhSDFWOWlsMlpPQkoyZFdMY0xCeXRVRmJvVXVpY2pYZ3puSGJaYkJrZWFqZ0lvYkQlMkJYQ1Z1OGF6em44b29pJTJGdXE3RmdKbDdMJTJCUXpkczl0bGx0NnFpdmVYdGNIdWFQbktWSjFBJTNEJlgtQW16LVNpZ25lZEhlYWRlcnM9aG9zdCZyZXNwb25zZS1leHBpcmVzPVNhdCUyQyUyMDE0JTIwSnVuJTIwMjAyNSUyMDAwJTNBMDAlM0EzNyUyMEdNVCZYLUFtei1TaWduYXR1cmU9YWVjYzBmYzJlYTlhNGNmZjQ1MTc1NGM0OWM3MmQzYzQwMDI4NzBlNjk3OTdjNjA1OTQyOTQyMTgyNmMzODExOA==")
$ws.Range("L2").Value = [System.Text.Encoding]::UTF8.GetString($l2bytes)

# Row 3: updated event - Harsh Brake / DAVID SERRANO (new timestamp & coords)
$ws.Range("A3").Value = "281474991395097-1749733783071"
$ws.Range("B3").Value = "Harsh Brake"
$ws.Range("C3").Value = "2025-06-12T07:09:43.071"
$ws.Range("D3").Value = "281474991395097"
$ws.Range("E3").Value = "125"
$ws.Range("F3").Value = "51834055"
$ws.Range("G3").Value = "DAVID SERRANO"
$ws.Range("H3").Value = 20.65073737
$ws.Range("I3").Value = -103.35702727
$ws.Range("J3").Value = 0.7843903303146362
$ws.Range("K3").Value = "No video URL"
$ws.Range("L3").Value = "No video URL"

# Restore default (Normal) style on the ID columns so no residual
# text-formatting style is left applied to the cells.
$idCols.Style = "Normal"
